$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D4:D9").FormulaR1C1 = '=IF(ISERROR(RC2),"ERROR",IF(ISERROR(RC3),"FAIL",IF(RC2=RC3,"PASS","FAIL")))'
$ws.Range("D3").Formula = '=IF(ISERROR(B3),"ERROR",IF(ISERROR(C3),"FAIL",IF(B3=C3,"PASS","FAIL")))'
